# Updates cryptos list price/volume columns (D, E) for rows 2-51 on Sheet1,
# matching the refreshed coinranking.com snapshot from the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.955.11'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').Value = '2.513.50'
$ws.Range('E3').Value = '  +0.49%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '530.49'
$ws.Range('E5').Value = '  -1.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.77'
$ws.Range('E6').Value = '  -3.47%  '
$ws.Range('E7').Value = '  +0.31%  '
$ws.Range('E8').Value = '  -1.52%  '
$ws.Range('D9').Value = '2.516.83'
$ws.Range('E9').Value = '  -0.25%  '
$ws.Range('E10').Value = '  -0.40%  '
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.44'
$ws.Range('E12').Value = '  -2.61%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.355'
$ws.Range('E13').Value = '  -0.01%  '
$ws.Range('D14').Value = '2.959.60'
$ws.Range('E14').Value = '  +0.47%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.05'
$ws.Range('E15').Value = '  -2.18%  '
$ws.Range('D16').Value = '58.922.52'
$ws.Range('E16').Value = '  +0.15%  '
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('D18').Value = '2.510.60'
$ws.Range('E18').Value = '  -0.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.00'
$ws.Range('E19').Value = '  -1.83%  '
$ws.Range('E20').Value = '  -0.78%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '321.82'
$ws.Range('E21').Value = '  -0.96%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('E23').Value = '  +0.67%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '62.32'
$ws.Range('E24').Value = '  +0.53%  '
$ws.Range('E25').Value = '  -3.83%  '
$ws.Range('E26').Value = '  +2.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.996'
$ws.Range('E27').Value = '  -0.36%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.80'
$ws.Range('E28').Value = '  +0.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.74'
$ws.Range('E29').Value = '  +1.06%  '
$ws.Range('D30').Value = '0.0₃0768'
$ws.Range('E30').Value = '  -0.66%  '
$ws.Range('E31').Value = '  -1.70%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '161.94'
$ws.Range('E32').Value = '  +3.62%  '
$ws.Range('E33').Value = '  +0.30%  '
$ws.Range('E34').Value = '  -6.16%  '
$ws.Range('E35').Value = '  -0.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.41'
$ws.Range('E36').Value = '  -1.36%  '
$ws.Range('E37').Value = '  -2.97%  '
$ws.Range('E38').Value = '  -1.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.00'
$ws.Range('E39').Value = '  +0.25%  '
$ws.Range('E40').Value = '  -1.67%  '
$ws.Range('E41').Value = '  -2.20%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.20'
$ws.Range('E42').Value = '  -8.55%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '279.26'
$ws.Range('E43').Value = '  -5.37%  '
$ws.Range('E45').Value = '  +0.97%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.595'
$ws.Range('E46').Value = '  -0.97%  '
$ws.Range('E47').Value = '  +0.04%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '121.46'
$ws.Range('E48').Value = '  -1.66%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '18.36'
$ws.Range('E49').Value = '  -1.13%  '
$ws.Range('E50').Value = '  -0.87%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0222'
$ws.Range('E51').Value = '  -2.30%  '
